$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.073.06'
$ws.Range('E2').Value = '  -0.97%  '
$ws.Range('D3').Value = '1.642.73'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('D5').Value = '217.38'
$ws.Range('E5').Value = '  -0.80%  '
$ws.Range('D6').Value = '0.5176'
$ws.Range('E6').Value = '  -3.25%  '
$ws.Range('E7').Value = '  -0.49%  '
$ws.Range('D8').Value = '0.2608'
$ws.Range('E8').Value = '  -1.93%  '
$ws.Range('D9').Value = '0.06275'
$ws.Range('D10').Value = '20.32'
$ws.Range('E10').Value = '  -1.99%  '
$ws.Range('D11').Value = '0.07758'
$ws.Range('E11').Value = '  -1.19%  '
$ws.Range('D12').Value = '4.465'
$ws.Range('E12').Value = '  -2.18%  '
$ws.Range('D13').Value = '1.626.22'
$ws.Range('E13').Value = '  -2.03%  '
$ws.Range('D14').Value = '1.869.94'
$ws.Range('D15').Value = '0.5547'
$ws.Range('E15').Value = '  +0.18%  '
$ws.Range('D16').Value = '0.0₅7978'
$ws.Range('E16').Value = '  -2.62%  '
$ws.Range('D17').Value = '64.63'
$ws.Range('E17').Value = '  -1.79%  '
$ws.Range('D18').Value = '26.073.14'
$ws.Range('E18').Value = '  -1.03%  '
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').Value = '4.608'
$ws.Range('E20').Value = '  -1.56%  '
$ws.Range('D21').Value = '192.60'
$ws.Range('E21').Value = '  -0.49%  '
$ws.Range('D22').Value = '10.05'
$ws.Range('E22').Value = '  -2.24%  '
$ws.Range('D23').Value = '5.933'
$ws.Range('E23').Value = '  -1.76%  '
$ws.Range('D24').Value = '1.006'
$ws.Range('E24').Value = '  -0.59%  '
$ws.Range('D25').Value = '146.84'
$ws.Range('E25').Value = '  +0.25%  '
$ws.Range('D26').Value = '0.1202'
$ws.Range('E26').Value = '  -2.49%  '
$ws.Range('D27').Value = '7.150'
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('D28').Value = '15.86'
$ws.Range('E28').Value = '  -1.66%  '
$ws.Range('D29').Value = '1.478'
$ws.Range('E29').Value = '  -1.53%  '
$ws.Range('D30').Value = '0.05629'
$ws.Range('E30').Value = '  -3.98%  '
$ws.Range('D31').Value = '1.267'
$ws.Range('E31').Value = '  -1.32%  '
$ws.Range('D32').Value = '3.464'
$ws.Range('E32').Value = '  -5.14%  '
$ws.Range('D33').Value = '3.359'
$ws.Range('E33').Value = '  +2.36%  '
$ws.Range('D34').Value = '1.585'
$ws.Range('E34').Value = '  -1.61%  '
$ws.Range('E35').Value = '  -1.38%  '
$ws.Range('D36').Value = '2.407'
$ws.Range('E36').Value = '  -0.47%  '
$ws.Range('D37').Value = '0.9343'
$ws.Range('E37').Value = '  -3.46%  '
$ws.Range('D38').Value = '0.5638'
$ws.Range('E38').Value = '  -3.26%  '
$ws.Range('D39').Value = '5.920'
$ws.Range('E39').Value = '  +1.05%  '
$ws.Range('D40').Value = '0.01573'
$ws.Range('E40').Value = '  -1.99%  '
$ws.Range('D41').Value = '1.057.75'
$ws.Range('E41').Value = '  +0.47%  '
$ws.Range('D42').Value = '1.005'
$ws.Range('E42').Value = '  -0.57%  '
$ws.Range('D43').Value = '0.8383'
$ws.Range('E43').Value = '  -3.56%  '
$ws.Range('D44').Value = '102.60'
$ws.Range('E44').Value = '  -2.55%  '
$ws.Range('D45').Value = '1.781.67'
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('D46').Value = '56.88'
$ws.Range('E46').Value = '  -1.61%  '
$ws.Range('D47').Value = '0.0₈105'
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').Value = '1.006'
$ws.Range('E48').Value = '  -0.89%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.05334'
$ws.Range('E49').Value = '  +3.20%  '
$ws.Range('D50').Value = '0.4330'
$ws.Range('E50').Value = '  -1.33%  '
$ws.Range('D51').Value = '7.930'
$ws.Range('E51').Value = '  -0.79%  '
